$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.587.08"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.960.50"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.19"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.616"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("E7").Value = "  +5.15%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  +5.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0791"
$ws.Range("E10").Value = "  -5.78%  "
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.22"
$ws.Range("E12").Value = "  +6.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.841"
$ws.Range("E13").Value = "  +5.07%  "
$ws.Range("D14").Value = "2.249.31"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.32"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("E16").Value = "  +3.54%  "
$ws.Range("D17").Value = "1.958.82"
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").Value = "36.562.88"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.74"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.99"
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.08"
$ws.Range("E22").Value = "  +2.95%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  +6.03%  "
$ws.Range("E25").Value = "  +4.49%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.143"
$ws.Range("E26").Value = "  +8.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.22"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.76"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.46"
$ws.Range("E29").Value = "  +1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.22"
$ws.Range("E30").Value = "  +10.62%  "
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.76"
$ws.Range("E32").Value = "  +4.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0617"
$ws.Range("E33").Value = "  -1.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.42"
$ws.Range("E34").Value = "  +7.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.58"
$ws.Range("E35").Value = "  +20.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.29"
$ws.Range("E36").Value = "  +8.36%  "
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.53"
$ws.Range("E39").Value = "  -8.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0983"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.17"
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.93"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").Value = "1.371.45"
$ws.Range("E45").Value = "  +3.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.36"
$ws.Range("E46").Value = "  +3.46%  "
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "44.72"
$ws.Range("E50").Value = "  +2.46%  "
$ws.Range("D51").Value = "2.139.48"
$ws.Range("E51").Value = "  +1.37%  "
